$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 10 (Objetivos:) which incorrectly held the professor's name ---
$objetivosPt = "Propiciar ao aluno uma visão básica sobre os principais métodos de determinação teórica da estrutura eletrônica dos materiais, com enfoque em sólidos cristalinos, mas também em materiais bidimensionais e nanoestruturados.O principal método de cálculo a ser empregado no curso será a Teoria do Funcional da Densidade(Density Functional Theory, DFT), em algumas de suas muitas variantes. Ao final do curso, o aluno estará apto a determinar propriedades dos materiais como estruturas de bandas, densidades de estados, superfícies de Fermi e constantes elásticas, usando um ou mais dos métodos e códigos computacionais apresentados em aula."
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# --- Insert a new row at 13 to hold the "Docentes responsáveis:" name, ---
# --- shifting the remaining rows (old 13..23) down to (14..24)        ---
$ws.Rows(13).Insert()
$ws.Range("A13").Clear()

# Match the body-cell formatting (style carries from the row-insert oddly) by
# copying formats from an existing B/C "data" cell before setting the value.
$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$docente = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# --- Row 14 ("Programa resumido:") content fix ---
$programaResumido = "Revisão de mecânica quântica; Revisão de física do estado sólido; Método de Hartree-Fock; Teoria do funcional da densidade; Métodos de ondas planas e pseudo-potenciais; Códigos computacionais"
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# --- Row 16 ("Programa:") content fix ---
$programa = "Revisão de mecânica quânticao Equação de Schrödingero Átomo do hidrogênio e orbitais atômicoso Notação de Diraco Princípio variacionalo Combinação linear de orbitais atômicosRevisão de física do estado sólidoo Espaço direto e recíprocoo Teorema de Blocho Zona de Brillouino Bandas de energia e densidade de estadoso Energia de Fermi e superfície de Fermio Aproximação de elétrons livresMétodo de Hartree-Focko Determinantes de Slatero Equação de Hartree-Focko Potencial de troca e correlaçãoo Algoritmo autoconsistenteTeoria do funcional da densidadeo Teoremas de Hohenberg-Kohno Equações de Kohn-Shamo Funcionais de troca e correlação: LDA, GGA, etc.Métodos de ondas planas e pseudo-potenciaiso Bases de ondas planaso Pseudo-potenciaiso Bases de ondas planas aumentadas e linearizadaso Método FP-LAPWCódigos computacionaiso Quantum Espressoo Elko Wien2ko VASP"
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Row 17 ("Syllabus:") was missing its English content; add it ---
$syllabus = "• Review of quantum mechanics: Schrödinger's equation; Hydrogen atom and atomic orbitals; Dirac notation; Variational principle; Linear combination of atomic orbitals. • Solid state physics review: Direct and reciprocal space; Bloch's Theorem; Brillouin zone; Energy bands and density of states; Fermi energy and Fermi surface; Free electrons Approximation. • Hartree-Fock method: Slater determinants; Hartree-Fock equation; Exchange and correlation potential; Self-consistent algorithm. • Density functional theory: Hohenberg-Kohn theorems; Kohn-Sham equations; Exchange and correlation functionals: LDA, GGA, etc. • Plane and pseudopotential wave methods: Plane wave bases; Pseudo-potentials; • Augmented and linearized plane wave bases: FP-LAPW method. • Computer codes: NWCHEM, Quantum Espresso, , Wien2k, exciting, VASP, etc."
$ws.Range("B17").Value = $syllabus
$ws.Range("C17").Value = $syllabus

# --- Row 19 ("Método:") content fix ---
$metodo = "Aulas expositivas, trabalhos e exercícios comentados."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Row 20 ("Critério:") content fix ---
$criterio = "Média aritmética de trabalhos propostos ao longo do curso."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Row 21 ("Norma de recuperação:") content fix ---
$norma = "Não haverá exame de recuperação"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# --- Row 22 ("Bibliografia:") was missing its content; add it ---
$bibliografia = @"
GRIFFITHS, D. J., Mecânica Quântica, Pearson.
ASHCROFT, N. W. Solid State Physics, Saunders College.
KITTEL, C. Introduction to Solid State Physics. John Wiley & Sons.
SUTTON, A. P. Electronic Structure of Materials, Oxford.
MORGON, N. H. e COUTINHO, K. (eds), Métodos de Química teórica e modelagem molecular, Livraria da Física
Editora.
VIANNA, J. D. M., FAZZIO, A., CANUTO, S., Teoria Quântica de moléculas e sólidos, Livraria da Física Editora.
COTTENIER, S. Density Functional Theory and the Family of (L)APW-methods: a step-by-step introduction
(apostila, disponível online)
THIJSSEN, J. M. Computational Physics, Cambridge.
TADMOR, E. B., MILLER, R. E. Modeling Materials  Continuum, atomistic and multiscale techniques,
Cambridge.
"@
$bibliografia = $bibliografia.TrimEnd("`r", "`n")
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# --- Column layout: column A no longer shares its width definition with column B ---
$ws.Columns("A").ColumnWidth = 29.8776041666667
